$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate existing data row 2 into new row 4
$ws.Range("A4").Value = "f3 add"
$ws.Range("B4").Value = "f3 city"
$ws.Range("C4").Value = "f3 first"
$ws.Range("D4").Value = "f3 last"
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = "'0303030"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = " "

# Duplicate existing data row 3 into new row 5
$ws.Range("A5").Value = "f2 arr"
$ws.Range("B5").Value = "f2 city"
$ws.Range("C5").Value = "{{address}} first name"
$ws.Range("D5").Value = "f2 last name"
$ws.Range("E5").Value = "{{address}}"
$ws.Range("F5").Value = "f2 222"
$ws.Range("G5").Value = "f2 zipcode"
